$d = $word.ActiveDocument

$old = "financial POA by my grandmother to help sell real estate, manage banking, and pay bills       July"
$new = "POA for my grandmother to help sell real-estate, manage banking, and pay bills                      July"

$range = $d.Content
$found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
Write-Output $found
